# Project DesignFirst save: update the "Integer min" rule value for R20
# (Rules sheet, cell C10) from 18 to 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 100.0
